# Insert two new weekly price rows for "Frutilla" at Feria Lagunitas de Puerto
# Montt (rows 46-47), pushing the existing data down by two rows.
#
# New dimension becomes A1:T113 (was A1:T111).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 46 (rows 46 and 47 in the old sheet, plus
# everything after, shift down by 2).
$ws.Range("46:47").Insert()

# Row 46: new "Especial" quality entry for 2021-10-26.
$row46 = @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44495, 10, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Especial", 300, 12000, 12000, 12000, "`$/bandeja 7 kilos", "Provincia de Melipilla", 1714, 7)
for ($i = 0; $i -lt $row46.Length; $i++) {
    $ws.Cells.Item(46, $i + 1).Value = $row46[$i]
}

# Row 47: new "Primera" quality entry for 2021-10-26.
$row47 = @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44495, 10, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Primera", 700, 9500, 10000, 9750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 1393, 7)
for ($i = 0; $i -lt $row47.Length; $i++) {
    $ws.Cells.Item(47, $i + 1).Value = $row47[$i]
}
